$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F7").Value = 5
$ws.Range("F10").Value = -2
$ws.Range("F13").Value = 1
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = -1
$ws.Range("F19").Value = 4
$ws.Range("F23").Value = -3
$ws.Range("F27").Value = -1
$ws.Range("F29").Value = 2
$ws.Range("F31").Value = 2
$ws.Range("F35").Value = -4
$ws.Range("F36").Value = -1
$ws.Range("F37").Value = -1
$ws.Range("F39").Value = -4
$ws.Range("F40").Value = 2
$ws.Range("F43").Value = -2
$ws.Range("F45").Value = 1
$ws.Range("F47").Value = 4
$ws.Range("F48").Value = 0
$ws.Range("F49").Value = 1
$ws.Range("F50").Value = -2
$ws.Range("F52").Value = 2
$ws.Range("F53").Value = -2
$ws.Range("F56").Value = 1
$ws.Range("F57").Value = 5
$ws.Range("F58").Value = -2
$ws.Range("F59").Value = 1
$ws.Range("F60").Value = -2
$ws.Range("F61").Value = -1
$ws.Range("F62").Value = 5
$ws.Range("F63").Value = 2
$ws.Range("F65").Value = -3
$ws.Range("F66").Value = 1
$ws.Range("F70").Value = 0
